# Apply the "cash & deposit done" update to the 存款 (deposits) sheet.
#
# Before the edit, sheet4 (存款) only had 6 columns (A..F) and row 1 was an
# accidental duplicate of row 2's data instead of being a header row. This
# change turns row 1 into a real header row and adds the same trailing
# metadata columns (property_category, category, date, legislator_name,
# legislator_id, source_file, index) that the other property sheets
# (土地/建物/汽車) already have.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Row 1: turn the old (accidental) data row into a proper header row ---
$ws.Cells.Item(1, 2).Value = "bank"            # B1
$ws.Cells.Item(1, 3).Value = "deposit_type"    # C1
$ws.Cells.Item(1, 4).Value = "currency"        # D1
$ws.Cells.Item(1, 5).Value = "owner"           # E1
$ws.Cells.Item(1, 6).Value = "total"           # F1
$ws.Cells.Item(1, 7).Value = "property_category"  # G1
$ws.Cells.Item(1, 8).Value = "category"           # H1
$ws.Cells.Item(1, 9).Value = "date"               # I1
$ws.Cells.Item(1, 10).Value = "legislator_name"   # J1
$ws.Cells.Item(1, 11).Value = "legislator_id"     # K1
$ws.Cells.Item(1, 12).Value = "source_file"       # L1
$ws.Cells.Item(1, 13).Value = "index"             # M1

# Copy the existing header style (bold + bordered, style index 1) onto the
# newly-added G1:M1 header cells so they match B1:F1.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2 (52: 第一商業銀行土城分行) ---
$ws.Cells.Item(2, 2).Value = "第一商業銀行土城分行"  # B2 bank
$ws.Cells.Item(2, 3).Value = "活期存款"              # C2 deposit_type
$ws.Cells.Item(2, 4).Value = "新臺幣"                # D2 currency
$ws.Cells.Item(2, 5).Value = "吳麗香"                # E2 owner
$ws.Cells.Item(2, 7).Value = "deposit"             # G2 property_category
$ws.Cells.Item(2, 8).Value = "normal"              # H2 category
$ws.Cells.Item(2, 9).Value = "2011-11-18"          # I2 date
$ws.Cells.Item(2, 10).Value = "盧嘉辰"              # J2 legislator_name
$ws.Cells.Item(2, 11).Value = 1715                 # K2 legislator_id
$ws.Cells.Item(2, 12).Value = "tmp94c1"            # L2 source_file
$ws.Cells.Item(2, 13).Value = 52                   # M2 index

# --- Row 3 (53: 中華郵政股份有限公司土城郵局) ---
$ws.Cells.Item(3, 3).Value = "活期存款"              # C3 deposit_type
$ws.Cells.Item(3, 4).Value = "新臺幣"                # D3 currency
$ws.Cells.Item(3, 7).Value = "deposit"             # G3 property_category
$ws.Cells.Item(3, 8).Value = "normal"              # H3 category
$ws.Cells.Item(3, 9).Value = "2011-11-18"          # I3 date
$ws.Cells.Item(3, 10).Value = "盧嘉辰"              # J3 legislator_name
$ws.Cells.Item(3, 11).Value = 1715                 # K3 legislator_id
$ws.Cells.Item(3, 12).Value = "tmp94c1"            # L3 source_file
$ws.Cells.Item(3, 13).Value = 53                   # M3 index

# --- Row 4 (54: 土藤會土城分行) ---
$ws.Cells.Item(4, 3).Value = "活期存款"              # C4 deposit_type
$ws.Cells.Item(4, 4).Value = "新臺幣"                # D4 currency
$ws.Cells.Item(4, 7).Value = "deposit"             # G4 property_category
$ws.Cells.Item(4, 8).Value = "normal"              # H4 category
$ws.Cells.Item(4, 9).Value = "2011-11-18"          # I4 date
$ws.Cells.Item(4, 10).Value = "盧嘉辰"              # J4 legislator_name
$ws.Cells.Item(4, 11).Value = 1715                 # K4 legislator_id
$ws.Cells.Item(4, 12).Value = "tmp94c1"            # L4 source_file
$ws.Cells.Item(4, 13).Value = 54                   # M4 index

# Copy the existing data-row style (style index 2) onto the newly-added
# G2:M4 cells so they match the rest of the data rows.
$ws.Range("B2").Copy()
$ws.Range("G2:M4").PasteSpecial(-4122)  # xlPasteFormats
